$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.089.54"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.57"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.99"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.06"
$ws.Range("E7").Value = "  +8.41%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0818"
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.76"
$ws.Range("E12").Value = "  +2.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.358.20"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.09"
$ws.Range("E14").Value = "  +3.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.759"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.057.45"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.057.45"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.79"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0830"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.79"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.58"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.96"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.29"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.48"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.05"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0604"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.29"
$ws.Range("E36").Value = "  +9.66%  "
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.536.68"
$ws.Range("E40").Value = "  +4.73%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0218"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.66"
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.65"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0930"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.03"
$ws.Range("E47").Value = "  -5.09%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.99"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.06"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.248.26"
$ws.Range("E51").Value = "  +1.15%  "
